$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)

$bvals = @(
    'Bitcoin',
    'Ethereum',
    'TetherUSD',
    'BNB',
    'Solana',
    'USDC',
    'XRP',
    'LidoStakedEther',
    'Toncoin',
    'Dogecoin',
    'Cardano',
    'TRON',
    'WrappedliquidstakedEther2.0',
    'WrappedBTC',
    'Avalanche',
    'ShibaInu',
    'WrappedEther',
    'Polkadot',
    'BitcoinCash',
    'Chainlink',
    'Uniswap',
    'Dai',
    'Litecoin',
    'Polygon',
    'WrappedeETH',
    'Kaspa',
    'Binance-PegBSC-USD',
    'PEPE',
    'InternetComputer(DFINITY)',
    'USDe',
    'EthereumClassic',
    'Monero',
    'PancakeSwap',
    'Aptos',
    'NEARProtocol',
    'ImmutableX',
    'SuiNetwork',
    'Stacks',
    'Fetch.AI',
    'OKB',
    'Filecoin',
    'Bittensor',
    'Stellar',
    'Mantle',
    'Hedera',
    'FirstDigitalUSD',
    'EnergySwap',
    'RenderToken',
    'VeChain',
    'WhiteBITCoin'
)

$cvals = @(
    'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc',
    'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth',
    'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt',
    'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb',
    'https://coinranking.com/coin/zNZHO_Sjf+solana-sol',
    'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc',
    'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp',
    'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth',
    'https://coinranking.com/coin/67YlI0K1b+toncoin-ton',
    'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge',
    'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada',
    'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx',
    'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth',
    'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc',
    'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax',
    'https://coinranking.com/coin/xz24e0BjL+shibainu-shib',
    'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',
    'https://coinranking.com/coin/25W7FG7om+polkadot-dot',
    'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',
    'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link',
    'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni',
    'https://coinranking.com/coin/MoTuySvg7+dai-dai',
    'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc',
    'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic',
    'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth',
    'https://coinranking.com/coin/V8GxkwWow+kaspa-kas',
    'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd',
    'https://coinranking.com/coin/03WI8NQPF+pepe-pepe',
    'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp',
    'https://coinranking.com/coin/exbfr2U-0+usde-usde',
    'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc',
    'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',
    'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake',
    'https://coinranking.com/coin/HGYj5JCv5+aptos-apt',
    'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',
    'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',
    'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui',
    'https://coinranking.com/coin/mMPrMcB7+stacks-stx',
    'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet',
    'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb',
    'https://coinranking.com/coin/ymQub4fuB+filecoin-fil',
    'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao',
    'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm',
    'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt',
    'https://coinranking.com/coin/jad286TjB+hedera-hbar',
    'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd',
    'https://coinranking.com/coin/SbWqqTui-+energyswap-ens',
    'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr',
    'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',
    'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
)

$dvals = @(
    '60.864.03',
    '2.619.63',
    '1.00',
    '513.95',
    '155.18',
    '0.997',
    '0.588',
    '2.634.84',
    '6.82',
    '0.105',
    '0.347',
    '0.129',
    '3.079.55',
    '60.830.11',
    '21.75',
    '0.0000141',
    '2.629.41',
    '4.75',
    '356.41',
    '10.62',
    '6.20',
    '1.00',
    '60.80',
    '0.424',
    '2.734.61',
    '0.167',
    '0.996',
    '0.0₃0849',
    '7.39',
    '0.999',
    '19.47',
    '152.64',
    '1.58',
    '5.94',
    '4.01',
    '1.20',
    '0.878',
    '1.49',
    '0.850',
    '36.38',
    '3.78',
    '294.89',
    '0.102',
    '0.625',
    '0.0558',
    '0.997',
    '19.90',
    '4.94',
    '0.0235',
    '10.31'
)

$evals = @(
    '  +0.67%  ',
    '  -0.12%  ',
    '  -0.09%  ',
    '  +1.10%  ',
    '  -1.55%  ',
    '  +0.39%  ',
    '  -0.29%  ',
    '  -1.02%  ',
    '  +5.03%  ',
    '  -0.36%  ',
    '  +0.85%  ',
    '  +1.87%  ',
    '  +0.25%  ',
    '  +0.52%  ',
    '  +0.03%  ',
    '  +0.02%  ',
    '  -0.87%  ',
    '  -0.97%  ',
    '  +2.83%  ',
    '  +0.93%  ',
    '  -0.01%  ',
    '  +0.27%  ',
    '  +1.23%  ',
    '  -0.07%  ',
    '  -0.33%  ',
    '  -0.07%  ',
    '  +0.57%  ',
    '  -1.74%  ',
    '  -2.83%  ',
    '  +0.16%  ',
    '  -0.32%  ',
    '  -2.46%  ',
    '  +0.52%  ',
    '  +2.73%  ',
    '  -1.12%  ',
    '  -1.24%  ',
    '  +4.43%  ',
    '  -0.06%  ',
    '  -0.55%  ',
    '  +2.71%  ',
    '  +0.13%  ',
    '  -4.54%  ',
    '  +1.17%  ',
    '  -1.74%  ',
    '  -3.02%  ',
    '  +0.52%  ',
    '  -1.27%  ',
    '  +0.28%  ',
    '  -0.62%  ',
    '  +0.30%  '
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
}